# Commit: "table update, removal of findrisk"
#
# The "Lindström et. al, 2013 (FINDRISC)" row in the predictive-models
# table on Sheet1 is removed entirely; every row below it shifts up by
# one. No other row content changes. Sheet2 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 held the FINDRISC study (Name/Year/Country/Ethnic comp/Sex comp/
# Variables/Outcome/Equation). Deleting the whole row shifts rows 5:11 up
# to 4:10, matching the diff exactly (dimension A1:K11 -> A1:K10).
$ws.Rows.Item(4).Delete()

# Leftover selection cursor after the edit, as recorded in the saved file.
$ws.Range("K5").Select()
